$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.284.68"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "1.900.97"
$ws.Range("E3").Value = "  +0.46%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "308.21"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5217"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.68%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3777"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07291"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.99%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "21.25"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.64%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.9027"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08211"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +7.28%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "96.93"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.71%  "

$ws.Range("D14").Value = "1.903.78"
$ws.Range("E14").Value = "  -1.32%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "5.300"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("E16").Value = "  +0.09%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000008660"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.73%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "14.59"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").Value = "27.311.11"
$ws.Range("E20").Value = "  +0.57%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.097"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.72"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +1.18%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.438"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.306"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.50%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "147.35"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "18.26"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.743"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "115.48"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.03%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "4.840"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.924"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.95%  "

$ws.Range("E31").Value = "  +0.35%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.05071"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.7990"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +2.75%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.235"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "3.444"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +4.81%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.948"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.29%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.594"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.48%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5706"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.28%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.02008"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "

$ws.Range("E40").Value = "  +0.33%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "9.016"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "6.575"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "116.26"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.83%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.1520"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.58%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.4887"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("E47").Value = "  +0.07%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.624"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "38.04"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.12%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "63.88"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05951"
$cell.NumberFormat = "General"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
